$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value as literal text (shared string), avoiding Excel's
# automatic conversion of date- or number-looking strings into dates/numbers.
function Set-TextValue($range, [string]$text) {
    $escaped = $text.Replace('"', '""')
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
    $excel.CutCopyMode = $false
}

# Update the report title (row 1) with the new date range
Set-TextValue $ws.Range("C1") "Facturacion de 2014-03-02 al 2014-03-04"

# Update the first invoice row (row 3) with the new values
$ws.Range("A3").Value = 1
Set-TextValue $ws.Range("B3") "2014-03-03"
Set-TextValue $ws.Range("C3") "15:29:57"
Set-TextValue $ws.Range("D3") "C"
Set-TextValue $ws.Range("E3") "100.00"

# Insert a new invoice row (pushes the totals row down, inherits row 3 style)
$ws.Rows("4").Insert()

$ws.Range("A4").Value = 2
Set-TextValue $ws.Range("B4") "2014-03-03"
Set-TextValue $ws.Range("C4") "16:37:52"
Set-TextValue $ws.Range("D4") "C"
Set-TextValue $ws.Range("E4") "167.00"

# Update the totals row (now row 6 after the insert)
Set-TextValue $ws.Range("A6") "Total Facturado"
Set-TextValue $ws.Range("E6") "267.00"
